$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 8 (shifts rows 9..192 up to 8..191, recalculating the
# sheet dimension and shared-string usage accordingly).
$ws.Rows(8).Delete()

# Restore the view/selection state recorded for the sheet after the edit.
$ws.Application.ActiveWindow.ScrollColumn = 15   # column O
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A8:XFD8").Select()
